# Case Header Report: add new "Resolution_Code" column to the "Data" table,
# following the existing "Responsible_Owner" column (diff: new table column,
# new shared string, worksheet header cell O1, table/autofilter range A1:O2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tbl = $ws.ListObjects.Item("Data")

# Appending a ListColumn extends the table (and its AutoFilter) by one
# column and grows the sheet dimension accordingly.
$newCol = $tbl.ListColumns.Add()
$hdr = $newCol.Range.Cells.Item(1, 1)
$hdr.Value = "Resolution_Code"

# Match the formatting Excel gives every other header cell in the table
# (bold header font on the accent fill) by copying it from a neighboring
# header cell instead of re-deriving the style by hand.
$ws.Range("N1").Copy() | Out-Null
$hdr.PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Return the active cell to A2, matching the saved selection state.
$ws.Range("A2").Select() | Out-Null
